# Windows 10 April 2018 Update SDK (17134) (#10)
#
# The "compatible with" blurb near the top of the Readme changes from:
#   "This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)"
# to:
#   "This sample is compatible with the Windows 10 April 2018 Update SDK (17134)"
#
# Apply the change as a series of small, unambiguous Find/Replace passes so
# each one only ever touches the single run it targets (each search string
# below is unique in the document).  Order matters: do the substitutions
# that don't interact with each other's search text first, then peel the
# "Fall "/"Creators" pair apart last.

$d = $word.ActiveDocument

# "...Update SDK (16299)" -> "...Update SDK (17134)"
$d.Content.Find.Execute("16299", $true, $false, $false, $false, $false, $true, 1, $false, "17134", 2) | Out-Null

# Drop the old " Update SDK (" run entirely; its text is reinserted (without
# the leading space) in front of "Creators" below.
$d.Content.Find.Execute(" Update SDK (", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# "Creators" -> "Update SDK ("
$d.Content.Find.Execute("Creators", $true, $false, $false, $false, $false, $true, 1, $false, "Update SDK (", 2) | Out-Null

# "Fall " -> "April 2018 "
$d.Content.Find.Execute("Fall ", $true, $false, $false, $false, $false, $true, 1, $false, "April 2018 ", 2) | Out-Null
